$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cntf"
$ws.Range("C2").Value = "Cntfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7676553333333332
$ws.Range("H2").Value = 2.302966
$ws.Range("I2").Value = 0.3736977786965754
$ws.Range("J2").Value = 0.3736977786965754
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009559333333333333
$ws.Range("N2").Value = 0.028678
$ws.Range("O2").Value = 0.001091248180460288
$ws.Range("P2").Value = 0.001091248180460288
$ws.Range("Q2").Value = 0.007338273216444443
$ws.Range("R2").Value = 0.06604445894799998
$ws.Range("S2").Value = 0.0004077970210446892
$ws.Range("T2").Value = 0.0004077970210446892

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cntf"
$ws.Range("C3").Value = "Cntfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7676553333333332
$ws.Range("H3").Value = 2.302966
$ws.Range("I3").Value = 0.3736977786965754
$ws.Range("J3").Value = 0.3736977786965754
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.501068666666667
$ws.Range("N3").Value = 25.503206
$ws.Range("O3").Value = 0.9704417024689271
$ws.Range("P3").Value = 0.9704417024689271
$ws.Range("Q3").Value = 6.525890700999555
$ws.Range("R3").Value = 58.73301630899599
$ws.Range("S3").Value = 0.3626519085671611
$ws.Range("T3").Value = 0.362651908567161

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cntf"
$ws.Range("C4").Value = "Cntfr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7676553333333332
$ws.Range("H4").Value = 2.302966
$ws.Range("I4").Value = 0.3736977786965754
$ws.Range("J4").Value = 0.3736977786965754
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1915436666666667
$ws.Range("N4").Value = 0.574631
$ws.Range("O4").Value = 0.02186571703696477
$ws.Range("P4").Value = 0.02186571703696477
$ws.Range("Q4").Value = 0.1470395172828889
$ws.Range("R4").Value = 1.323355655546
$ws.Range("S4").Value = 0.008171169886321598
$ws.Range("T4").Value = 0.008171169886321598

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cntf"
$ws.Range("C5").Value = "Cntfr"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7676553333333332
$ws.Range("H5").Value = 2.302966
$ws.Range("I5").Value = 0.3736977786965754
$ws.Range("J5").Value = 0.3736977786965754
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05782766666666667
$ws.Range("N5").Value = 0.173483
$ws.Range("O5").Value = 0.006601332313647817
$ws.Range("P5").Value = 0.006601332313647817
$ws.Range("Q5").Value = 0.04439171673088888
$ws.Range("R5").Value = 0.3995254505779999
$ws.Range("S5").Value = 0.002466903222048114
$ws.Range("T5").Value = 0.002466903222048114

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cntf"
$ws.Range("C6").Value = "Cntfr"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.811404
$ws.Range("H6").Value = 2.434212
$ws.Range("I6").Value = 0.3949948098567449
$ws.Range("J6").Value = 0.3949948098567449
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.009559333333333333
$ws.Range("N6").Value = 0.028678
$ws.Range("O6").Value = 0.001091248180460288
$ws.Range("P6").Value = 0.001091248180460288
$ws.Range("Q6").Value = 0.007756481304
$ws.Range("R6").Value = 0.069808331736
$ws.Range("S6").Value = 0.0004310373675474302
$ws.Range("T6").Value = 0.0004310373675474302

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cntf"
$ws.Range("C7").Value = "Cntfr"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.811404
$ws.Range("H7").Value = 2.434212
$ws.Range("I7").Value = 0.3949948098567449
$ws.Range("J7").Value = 0.3949948098567449
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.501068666666667
$ws.Range("N7").Value = 25.503206
$ws.Range("O7").Value = 0.9704417024689271
$ws.Range("P7").Value = 0.9704417024689271
$ws.Range("Q7").Value = 6.897801120408
$ws.Range("R7").Value = 62.080210083672
$ws.Range("S7").Value = 0.3833194357437697
$ws.Range("T7").Value = 0.3833194357437697

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cntf"
$ws.Range("C8").Value = "Cntfr"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.811404
$ws.Range("H8").Value = 2.434212
$ws.Range("I8").Value = 0.3949948098567449
$ws.Range("J8").Value = 0.3949948098567449
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1915436666666667
$ws.Range("N8").Value = 0.574631
$ws.Range("O8").Value = 0.02186571703696477
$ws.Range("P8").Value = 0.02186571703696477
$ws.Range("Q8").Value = 0.155419297308
$ws.Range("R8").Value = 1.398773675772
$ws.Range("S8").Value = 0.008636844743397286
$ws.Range("T8").Value = 0.008636844743397286

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cntf"
$ws.Range("C9").Value = "Cntfr"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.811404
$ws.Range("H9").Value = 2.434212
$ws.Range("I9").Value = 0.3949948098567449
$ws.Range("J9").Value = 0.3949948098567449
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.05782766666666667
$ws.Range("N9").Value = 0.173483
$ws.Range("O9").Value = 0.006601332313647817
$ws.Range("P9").Value = 0.006601332313647817
$ws.Range("Q9").Value = 0.046921600044
$ws.Range("R9").Value = 0.422294400396
$ws.Range("S9").Value = 0.002607492002030505
$ws.Range("T9").Value = 0.002607492002030505

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Cntf"
$ws.Range("C10").Value = "Cntfr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.475155
$ws.Range("H10").Value = 1.425465
$ws.Range("I10").Value = 0.2313074114466796
$ws.Range("J10").Value = 0.2313074114466796
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.009559333333333333
$ws.Range("N10").Value = 0.028678
$ws.Range("O10").Value = 0.001091248180460288
$ws.Range("P10").Value = 0.001091248180460288
$ws.Range("Q10").Value = 0.00454216503
$ws.Range("R10").Value = 0.04087948526999999
$ws.Range("S10").Value = 0.0002524137918681683
$ws.Range("T10").Value = 0.0002524137918681683

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Cntf"
$ws.Range("C11").Value = "Cntfr"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.475155
$ws.Range("H11").Value = 1.425465
$ws.Range("I11").Value = 0.2313074114466796
$ws.Range("J11").Value = 0.2313074114466796
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.501068666666667
$ws.Range("N11").Value = 25.503206
$ws.Range("O11").Value = 0.9704417024689271
$ws.Range("P11").Value = 0.9704417024689271
$ws.Range("Q11").Value = 4.03932528231
$ws.Range("R11").Value = 36.35392754079
$ws.Range("S11").Value = 0.2244703581579963
$ws.Range("T11").Value = 0.2244703581579963

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Cntf"
$ws.Range("C12").Value = "Cntfr"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.475155
$ws.Range("H12").Value = 1.425465
$ws.Range("I12").Value = 0.2313074114466796
$ws.Range("J12").Value = 0.2313074114466796
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.1915436666666667
$ws.Range("N12").Value = 0.574631
$ws.Range("O12").Value = 0.02186571703696477
$ws.Range("P12").Value = 0.02186571703696477
$ws.Range("Q12").Value = 0.091012930935
$ws.Range("R12").Value = 0.819116378415
$ws.Range("S12").Value = 0.005057702407245882
$ws.Range("T12").Value = 0.005057702407245882

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Cntf"
$ws.Range("C13").Value = "Cntfr"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.475155
$ws.Range("H13").Value = 1.425465
$ws.Range("I13").Value = 0.2313074114466796
$ws.Range("J13").Value = 0.2313074114466796
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05782766666666667
$ws.Range("N13").Value = 0.173483
$ws.Range("O13").Value = 0.006601332313647817
$ws.Range("P13").Value = 0.006601332313647817
$ws.Range("Q13").Value = 0.027477104955
$ws.Range("R13").Value = 0.247293944595
$ws.Range("S13").Value = 0.001526937089569197
$ws.Range("T13").Value = 0.001526937089569197
